$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.03148819180643
$ws.Range("C2").Value = 9.593674530799468
$ws.Range("E2").Value = 12.74524976387822
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.673189197126387
$ws.Range("K2").Value = 9.418277300808496
$ws.Range("L2").Value = 9.951288756319935
$ws.Range("M2").Value = 14.46702459954551
$ws.Range("N2").Value = 20.77650444022321
$ws.Range("O2").Value = 26.42444030374187

$ws.Range("B3").Value = 12.85150573726562
$ws.Range("C3").Value = 9.595493443088982
$ws.Range("E3").Value = 12.77179413609868
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.674819255093851
$ws.Range("K3").Value = 9.285767747287739
$ws.Range("L3").Value = 9.958310512857667
$ws.Range("M3").Value = 14.44385634529189
$ws.Range("N3").Value = 20.83781109935131
$ws.Range("O3").Value = 26.51693559186113

$ws.Range("B4").Value = 12.74216574062311
$ws.Range("C4").Value = 9.596946724846834
$ws.Range("E4").Value = 12.78982880706979
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.675873612738297
$ws.Range("K4").Value = 9.204943792927869
$ws.Range("L4").Value = 9.963924471670028
$ws.Range("M4").Value = 14.43167193604012
$ws.Range("N4").Value = 20.87725191444903
$ws.Range("O4").Value = 26.57835323992428

$ws.Range("B5").Value = 12.6979565348799
$ws.Range("C5").Value = 9.597623932879811
$ws.Range("E5").Value = 12.79761497541842
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.676316764708715
$ws.Range("K5").Value = 9.172180848905978
$ws.Range("L5").Value = 9.966540352302111
$ws.Range("M5").Value = 14.42722395710906
$ws.Range("N5").Value = 20.89377786695563
$ws.Range("O5").Value = 26.60454380150546

$ws.Range("B6").Value = 12.69063818917687
$ws.Range("C6").Value = 9.597741526611621
$ws.Range("E6").Value = 12.79893425786183
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.676391165881206
$ws.Range("K6").Value = 9.166752182994317
$ws.Range("L6").Value = 9.966994552982312
$ws.Range("M6").Value = 14.42651672725979
$ws.Range("N6").Value = 20.89654941758282
$ws.Range("O6").Value = 26.60896290929208

$ws.Range("B7").Value = 12.7415680414088
$ws.Range("C7").Value = 9.59695551331146
$ws.Range("E7").Value = 12.78993204473043
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.675879534554796
$ws.Range("K7").Value = 9.20450118743644
$ws.Range("L7").Value = 9.96395842094179
$ws.Range("M7").Value = 14.43160984953768
$ws.Range("N7").Value = 20.87747295128389
$ws.Range("O7").Value = 26.57870174948258

$ws.Range("B8").Value = 12.96921961848429
$ws.Range("C8").Value = 9.594232085279977
$ws.Range("E8").Value = 12.75404206619561
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.673740162732551
$ws.Range("K8").Value = 9.372499455456078
$ws.Range("L8").Value = 9.953439888942315
$ws.Range("M8").Value = 14.45861509150914
$ws.Range("N8").Value = 20.79727053445455
$ws.Range("O8").Value = 26.45537232969447

$ws.Range("B9").Value = 13.42256849938506
$ws.Range("C9").Value = 9.591544854090651
$ws.Range("E9").Value = 12.69742596350339
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.669967460214702
$ws.Range("K9").Value = 9.704529013272479
$ws.Range("L9").Value = 9.943120169614339
$ws.Range("M9").Value = 14.52757730701004
$ws.Range("N9").Value = 20.65420175177243
$ws.Range("O9").Value = 26.25025178827579

$ws.Range("B10").Value = 13.75672673170224
$ws.Range("C10").Value = 9.591166727569243
$ws.Range("E10").Value = 12.66420294012902
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.667450672999058
$ws.Range("K10").Value = 9.947824651618161
$ws.Range("L10").Value = 9.941781306677495
$ws.Range("M10").Value = 14.58773019552329
$ws.Range("N10").Value = 20.55766448206705
$ws.Range("O10").Value = 26.12198001054916

$ws.Range("B11").Value = 13.90832831704942
$ws.Range("C11").Value = 9.591336743720575
$ws.Range("E11").Value = 12.6509026860457
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.666360541002083
$ws.Range("K11").Value = 10.05790791418587
$ws.Range("L11").Value = 9.942517878757858
$ws.Range("M11").Value = 14.61709130505821
$ws.Range("N11").Value = 20.5155908662139
$ws.Range("O11").Value = 26.06850500996562

$ws.Range("B12").Value = 13.96562444674974
$ws.Range("C12").Value = 9.591449913287956
$ws.Range("E12").Value = 12.64612654230871
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.665955569920593
$ws.Range("K12").Value = 10.09947138362029
$ws.Range("L12").Value = 9.942989317790301
$ws.Range("M12").Value = 14.62849088113044
$ws.Range("N12").Value = 20.4999221564455
$ws.Range("O12").Value = 26.04895743976089

$ws.Range("B13").Value = 13.95329051724809
$ws.Range("C13").Value = 9.591423376650399
$ws.Range("E13").Value = 12.64714359591745
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.666042439669529
$ws.Range("K13").Value = 10.09052598437046
$ws.Range("L13").Value = 9.942879239268477
$ws.Range("M13").Value = 14.62602337296629
$ws.Range("N13").Value = 20.50328498588122
$ws.Range("O13").Value = 26.05313610974887

$ws.Range("B14").Value = 13.91304465037991
$ws.Range("C14").Value = 9.591345078733683
$ws.Range("E14").Value = 12.65050453362105
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.666327066909591
$ws.Range("K14").Value = 10.06133003398145
$ws.Range("L14").Value = 9.942552811944552
$ws.Range("M14").Value = 14.61802355577238
$ws.Range("N14").Value = 20.51429651467194
$ws.Range("O14").Value = 26.06688273960858

$ws.Range("B15").Value = 13.88837674130752
$ws.Range("C15").Value = 9.591303460887962
$ws.Range("E15").Value = 12.65259710266145
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.666502428984089
$ws.Range("K15").Value = 10.04342960785456
$ws.Range("L15").Value = 9.942377906447641
$ws.Range("M15").Value = 14.61315986804988
$ws.Range("N15").Value = 20.52107569505816
$ws.Range("O15").Value = 26.07539443140935

$ws.Range("B16").Value = 13.74680639205192
$ws.Range("C16").Value = 9.591162465377842
$ws.Range("E16").Value = 12.66510859793769
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.667523014556783
$ws.Range("K16").Value = 9.940615292337146
$ws.Range("L16").Value = 9.941760164983641
$ws.Range("M16").Value = 14.58585105794499
$ws.Range("N16").Value = 20.56045105027766
$ws.Range("O16").Value = 26.12557294218729

$ws.Range("B17").Value = 13.65981444716912
$ws.Range("C17").Value = 9.591163290456164
$ws.Range("E17").Value = 12.67324813061149
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.668163111353281
$ws.Range("K17").Value = 9.877363138338938
$ws.Range("L17").Value = 9.941725171424681
$ws.Range("M17").Value = 14.5696053256777
$ws.Range("N17").Value = 20.58507741225185
$ws.Range("O17").Value = 26.15760554427372

$ws.Range("B18").Value = 13.60974361204027
$ws.Range("C18").Value = 9.591195982093186
$ws.Range("E18").Value = 12.67810045132686
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.668536435247023
$ws.Range("K18").Value = 9.840928505700381
$ws.Range("L18").Value = 9.941831764518732
$ws.Range("M18").Value = 14.56044958169398
$ws.Range("N18").Value = 20.59941526071522
$ws.Range("O18").Value = 26.17648878599985

$ws.Range("B19").Value = 13.59278606419348
$ws.Range("C19").Value = 9.591212596224754
$ws.Range("E19").Value = 12.67977268758647
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.668663723186236
$ws.Range("K19").Value = 9.828584306938454
$ws.Range("L19").Value = 9.941889648238677
$ws.Range("M19").Value = 14.55738213649496
$ws.Range("N19").Value = 20.60429963027468
$ws.Range("O19").Value = 26.18296111870881

$ws.Range("B20").Value = 13.6690789402896
$ws.Range("C20").Value = 9.591159871136906
$ws.Range("E20").Value = 12.67236400144457
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.668094438482316
$ws.Range("K20").Value = 9.884102267123374
$ws.Range("L20").Value = 9.941715788352539
$ws.Range("M20").Value = 14.57131525872624
$ws.Range("N20").Value = 20.58243795457293
$ws.Range("O20").Value = 26.15414811348892

$ws.Range("B21").Value = 13.92486929625092
$ws.Range("C21").Value = 9.591366755795235
$ws.Range("E21").Value = 12.64951028116412
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.666243252557307
$ws.Range("K21").Value = 10.06990921325328
$ws.Range("L21").Value = 9.942643474972792
$ws.Range("M21").Value = 14.62036571492911
$ws.Range("N21").Value = 20.5110550155152
$ws.Range("O21").Value = 26.06282595418619

$ws.Range("B22").Value = 14.09136615120563
$ws.Range("C22").Value = 9.591786156386233
$ws.Range("E22").Value = 12.63609153044425
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.665079065807187
$ws.Range("K22").Value = 10.19061334679406
$ws.Range("L22").Value = 9.944371387966065
$ws.Range("M22").Value = 14.65405884534896
$ws.Range("N22").Value = 20.46593838701206
$ws.Range("O22").Value = 26.0072353729112

$ws.Range("B23").Value = 14.00258276623135
$ws.Range("C23").Value = 9.591536442213508
$ws.Range("E23").Value = 12.64311464241388
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.665696247560807
$ws.Range("K23").Value = 10.12627022150058
$ws.Range("L23").Value = 9.943346884906994
$ws.Range("M23").Value = 14.63592857857082
$ws.Range("N23").Value = 20.48987779906089
$ws.Range("O23").Value = 26.03653020092765

$ws.Range("B24").Value = 13.66489063954184
$ws.Range("C24").Value = 9.59116131663524
$ws.Range("E24").Value = 12.67276317798461
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.66812546890141
$ws.Range("K24").Value = 9.881055722562598
$ws.Range("L24").Value = 9.941719635676883
$ws.Range("M24").Value = 14.57054162395091
$ws.Range("N24").Value = 20.58363069314887
$ws.Range("O24").Value = 26.15570976249457

$ws.Range("B25").Value = 13.29951078482256
$ws.Range("C25").Value = 9.591989958028929
$ws.Range("E25").Value = 12.71127034798656
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.67094310508408
$ws.Range("K25").Value = 9.6146612361987
$ws.Range("L25").Value = 9.944812377711559
$ws.Range("M25").Value = 14.5072344671031
$ws.Range("N25").Value = 20.69139354216528
$ws.Range("O25").Value = 26.30180621392222
